$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73-106 down to 74-107
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record
$ws.Cells.Item(73, 1).Value = 5
$ws.Cells.Item(73, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(73, 3).Value = "Maule"
$ws.Cells.Item(73, 4).Value = 45141
$ws.Cells.Item(73, 5).Value = 7
$ws.Cells.Item(73, 6).Value = 100112040
$ws.Cells.Item(73, 7).Value = "Cilantro"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 150
$ws.Cells.Item(73, 11).Value = 8000
$ws.Cells.Item(73, 12).Value = 8000
$ws.Cells.Item(73, 13).Value = 8000
$ws.Cells.Item(73, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 222
$ws.Cells.Item(73, 17).Value = 36
$ws.Cells.Item(73, 18).Value = "Hortaliza"
